# "change name shipping option"
#
# Renames the three "shippingOption" choice values that are shared across
# the FPA* order sheets:
#   "Standard shipping"                            -> "Package delivery (extra costs)"
#   "Pickup at factory"                             -> "Pick-up at factory (no costs)"
#   "Special packaging / via freight forwarding"    -> "Freight delivery / sepcial packaging (extra costs)"
#
# The values live in column N (O on "FPA008-009", which has an extra
# leading column) on rows 2-4 of each of the four FPA sheets.

$wb = $excel.ActiveWorkbook

$newPackageDelivery = "Package delivery (extra costs)"
$newPickupAtFactory = "Pick-up at factory (no costs)"
$newFreightDelivery = "Freight delivery / sepcial packaging (extra costs)"

# FPA001 - shippingOption is column N
$ws1 = $wb.Worksheets.Item("FPA001")
$ws1.Range("N2").Value = $newPackageDelivery
$ws1.Range("N3").Value = $newPickupAtFactory
$ws1.Range("N4").Value = $newFreightDelivery

# FPA002-003-005-007 - shippingOption is column N
$ws2 = $wb.Worksheets.Item("FPA002-003-005-007")
$ws2.Range("N2").Value = $newPackageDelivery
$ws2.Range("N3").Value = $newPickupAtFactory
$ws2.Range("N4").Value = $newFreightDelivery

# FPA004-006-010 - shippingOption is column N
$ws3 = $wb.Worksheets.Item("FPA004-006-010")
$ws3.Range("N2").Value = $newPackageDelivery
$ws3.Range("N3").Value = $newPickupAtFactory
$ws3.Range("N4").Value = $newFreightDelivery

# FPA008-009 - shippingOption is column O (extra customerEmail column shifts things)
$ws4 = $wb.Worksheets.Item("FPA008-009")
$ws4.Range("O2").Value = $newPackageDelivery
$ws4.Range("O3").Value = $newPickupAtFactory
$ws4.Range("O4").Value = $newFreightDelivery

# Widen column N on "FPA002-003-005-007" so the longer new text is fully
# visible (author manually resized it after the rename).
$ws2.Columns.Item(14).ColumnWidth = 39.6

# Reflect the view-state changes the author left behind while verifying the
# edits: FPA001 becomes the active sheet/tab with N3 selected, the
# previously-active FPA002-003-005-007 sheet is left with J20 selected, and
# BTMI015 is left with M29 selected.
$ws7 = $wb.Worksheets.Item("BTMI015")
$ws7.Activate()
$ws7.Range("M29").Select()

$ws2.Activate()
$ws2.Range("J20").Select()

$ws1.Activate()
$ws1.Range("N3").Select()
